$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.550.21'
$ws.Range("E2").Value = '  +14.01%  '
$ws.Range("D3").Value = '1.826.29'
$ws.Range("E3").Value = '  +8.64%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''233.53'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.70%  '
$ws.Range("D6").Value = '''0.552'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.50%  '
$ws.Range("D7").Value = '''1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").Value = '''32.29'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +8.08%  '
$ws.Range("D9").Value = '''46.33'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.08%  '
$ws.Range("E10").Value = '  +8.42%  '
$ws.Range("D11").Value = '''0.0684'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +10.00%  '
$ws.Range("E12").Value = '  +3.21%  '
$ws.Range("D13").Value = '2.091.81'
$ws.Range("E13").Value = '  +8.95%  '
$ws.Range("D14").Value = '1.823.24'
$ws.Range("E14").Value = '  +8.26%  '
$ws.Range("D15").Value = '''0.647'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.63%  '
$ws.Range("D16").Value = '34.507.26'
$ws.Range("E16").Value = '  +13.84%  '
$ws.Range("D17").Value = '''10.43'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.17%  '
$ws.Range("D18").Value = '''4.33'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +8.26%  '
$ws.Range("D19").Value = '''71.38'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +8.53%  '
$ws.Range("D20").Value = '''263.44'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +6.46%  '
$ws.Range("D21").Value = '0.0₃0762'
$ws.Range("E21").Value = '  +6.07%  '
$ws.Range("E22").Value = '  -0.16%  '
$ws.Range("D23").Value = '''10.59'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.84%  '
$ws.Range("D24").Value = '''4.43'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.10%  '
$ws.Range("D25").Value = '''2.19'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.05%  '
$ws.Range("D26").Value = '''162.02'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.01%  '
$ws.Range("E27").Value = '  +7.52%  '
$ws.Range("E28").Value = '  +5.77%  '
$ws.Range("E29").Value = '  +6.70%  '
$ws.Range("E30").Value = '  -0.03%  '
$ws.Range("D31").Value = '''3.87'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +10.74%  '
$ws.Range("D32").Value = '''0.0519'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.66%  '
$ws.Range("D33").Value = '''1.21'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.64%  '
$ws.Range("E34").Value = '  +9.01%  '
$ws.Range("D35").Value = '1.597.93'
$ws.Range("E35").Value = '  +6.73%  '
$ws.Range("E36").Value = '  +7.03%  '
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").Value = '''1.07'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.05%  '
$ws.Range("B38").Value = 'Aave'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D38").Value = '''86.51'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +9.78%  '
$ws.Range("D39").Value = '''0.636'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +8.68%  '
$ws.Range("D40").Value = '''0.0189'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.90%  '
$ws.Range("D41").Value = '''2.82'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.71%  '
$ws.Range("B42").Value = 'HuobiToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D42").Value = '''2.37'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.41%  '
$ws.Range("B43").Value = 'ARBITRUM'
$ws.Range("C43").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D43").Value = '''0.927'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +8.56%  '
$ws.Range("E44").Value = '  +7.18%  '
$ws.Range("D45").Value = '''0.0529'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.27%  '
$ws.Range("D46").Value = '''1.07'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.82%  '
$ws.Range("D47").Value = '1.981.97'
$ws.Range("D48").Value = '''54.37'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.80%  '
$ws.Range("D49").Value = '''5.78'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +6.41%  '
$ws.Range("E50").Value = '  -0.09%  '
$ws.Range("D51").Value = '''11.51'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +24.60%  '
